$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Analysis of lysosomal enzyme activities in  induced pluripotent stem cell, neural progenitor  cell, and neuron models as potential biomarkers  of Huntington’s Disease"
$ws.Range("B12").Value = "['Jimbo']"
